$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-05-31 Saturday" "2025-06-01 Sunday"

Replace-Text "656÷3=" "874÷9="
Replace-Text "208÷5=" "342÷9="
Replace-Text "658÷3=" "438÷4="
Replace-Text "191÷7=" "268÷3="
Replace-Text "757÷2=" "565÷4="

Replace-Text "626÷4=" "624÷9="
Replace-Text "230÷8=" "593÷2="
Replace-Text "242÷3=" "703÷6="
Replace-Text "738÷3=" "593÷6="
Replace-Text "523÷6=" "997÷6="

Replace-Text "483÷5=" "226÷3="
Replace-Text "152÷4=" "763÷9="
Replace-Text "943÷6=" "658÷5="
Replace-Text "892÷2=" "454÷8="
Replace-Text "958÷5=" "118÷3="

Replace-Text "871÷5=" "521÷4="
Replace-Text "216÷3=" "664÷4="
Replace-Text "659÷4=" "973÷4="
Replace-Text "950÷4=" "939÷3="
Replace-Text "485÷3=" "957÷9="

Replace-Text "569÷2=" "696÷9="
Replace-Text "786÷2=" "654÷6="
Replace-Text "134÷5=" "769÷6="
Replace-Text "151÷5=" "921÷9="
Replace-Text "587÷6=" "512÷9="
